$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-30 Friday", "2024-08-31 Saturday"),
    @("57×20=1140", "93×25=2325"),
    @("62×58=3596", "34×92=3128"),
    @("44×34=1496", "61×27=1647"),
    @("97×68=6596", "62×85=5270"),
    @("72×56=4032", "57×47=2679"),
    @("67×31=2077", "94×64=6016"),
    @("22×32=704", "33×35=1155"),
    @("47×94=4418", "37×92=3404"),
    @("35×23=805", "82×94=7708"),
    @("81×64=5184", "47×52=2444"),
    @("12×52=624", "87×76=6612"),
    @("95×81=7695", "78×93=7254"),
    @("17×95=1615", "49×40=1960"),
    @("61×35=2135", "49×84=4116"),
    @("17×69=1173", "62×52=3224"),
    @("26×97=2522", "81×43=3483"),
    @("70×77=5390", "49×17=833"),
    @("19×98=1862", "92×27=2484"),
    @("98×63=6174", "19×51=969"),
    @("58×30=1740", "54×35=1890"),
    @("82×52=4264", "49×43=2107"),
    @("39×98=3822", "40×66=2640"),
    @("63×13=819", "38×88=3344"),
    @("90×84=7560", "50×36=1800"),
    @("66×35=2310", "68×28=1904")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
